$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "teste"
$ws.Range("B3").Value = "notaaa"
$ws.Range("A4").Value = 17188
$ws.Range("B4").Value = 6377

$ws.Range("B5").Select()
